# Horarios 141 — 31/12 14:57 scrape update (LP1912 + LP1912-215/6203-6173 feeds)
# Appends newly scraped arrival rows to the three sheets and refreshes the
# "Ultima actualizacion" / "Total filas" header cells + used-range dimension.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"  (columns: A blank | B Hora_Scrap | C Hora_Llegada |
#                      D Linea | E Minutos | F Parada | G Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 11:57:24"
$ws1.Range("A3").Value = "Total filas: 918"

$sheet1Rows = @(
    @(903, "11:57:13", "11:59", "11_ETCHEVERRY",        2, "LP1912", "31/12/2025"),
    @(904, "11:57:13", "12:04", "23_HERNANDEZ",         7, "LP1912", "31/12/2025"),
    @(905, "11:57:13", "12:06", "17_ROMERO",             9, "LP1912", "31/12/2025"),
    @(906, "11:57:13", "12:10", "16_SANTA ANA",         13, "LP1912", "31/12/2025"),
    @(907, "11:57:13", "12:18", "10_OLMOS",             21, "LP1912", "31/12/2025"),
    @(908, "11:57:13", "12:18", "15_ABASTO",            21, "LP1912", "31/12/2025"),
    @(909, "11:57:13", "12:29", "215C_EL PATO",         32, "LP1912", "31/12/2025"),
    @(910, "11:57:13", "12:34", "23_HERNANDEZ",         37, "LP1912", "31/12/2025"),
    @(911, "11:57:13", "12:41", "15X38_ABASTO",         44, "LP1912", "31/12/2025"),
    @(912, "11:57:13", "12:52", "15_ABASTO",            55, "LP1912", "31/12/2025"),
    @(913, "11:57:13", "12:59", "16_SANTA ANA",         62, "LP1912", "31/12/2025"),
    @(914, "11:57:13", "13:01", "215C_EL PATO",         64, "LP1912", "31/12/2025"),
    @(915, "11:57:13", "13:05", "23_HERNANDEZ",         68, "LP1912", "31/12/2025"),
    @(916, "11:57:13", "13:07", "14_ABASTO",            70, "LP1912", "31/12/2025"),
    @(917, "11:57:13", "13:22", "17_ROMERO",            85, "LP1912", "31/12/2025"),
    @(918, "11:57:13", "13:30", "10_OLMOS",             93, "LP1912", "31/12/2025"),
    @(919, "11:57:13", "13:32", "16_P MOR-SANTA ANA",   95, "LP1912", "31/12/2025")
)

foreach ($r in $sheet1Rows) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 1).Style = "Normal"
    $ws1.Cells.Item($rowNum, 2).Value = $r[1]
    $ws1.Cells.Item($rowNum, 3).Value = $r[2]
    $ws1.Cells.Item($rowNum, 4).Value = $r[3]
    $ws1.Cells.Item($rowNum, 5).Value = $r[4]
    $ws1.Cells.Item($rowNum, 6).Value = $r[5]
    $ws1.Cells.Item($rowNum, 7).Value = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215"  (columns: A blank | B Fecha | C Hora_Scrap |
#                          D Hora_Llegada | E Linea | F Minutos | G Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 11:57:24"
$ws2.Range("A3").Value = "Total filas: 70"

$sheet2Rows = @(
    @(70, "31/12/2025", "11:57:13", "12:29", "215C_EL PATO", 32, "LP1912"),
    @(71, "31/12/2025", "11:57:13", "13:01", "215C_EL PATO", 64, "LP1912")
)

foreach ($r in $sheet2Rows) {
    $rowNum = $r[0]
    $ws2.Cells.Item($rowNum, 1).Style = "Normal"
    $ws2.Cells.Item($rowNum, 2).Value = $r[1]
    $ws2.Cells.Item($rowNum, 3).Value = $r[2]
    $ws2.Cells.Item($rowNum, 4).Value = $r[3]
    $ws2.Cells.Item($rowNum, 5).Value = $r[4]
    $ws2.Cells.Item($rowNum, 6).Value = $r[5]
    $ws2.Cells.Item($rowNum, 7).Value = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173"  (columns: A blank | B Fecha | C Hora_Scrap |
#                         D Hora_Llegada | E Linea | F Minutos | G Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 11:57:24"
$ws3.Range("A3").Value = "Total filas: 107"

$sheet3Rows = @(
    @(107, "31/12/2025", "11:57:24", "13:09", "215B_LP-P MOR-1 Y 57", 72, "L6173"),
    @(108, "31/12/2025", "11:57:24", "13:14", "215A_LA PLATA",        77, "L6173")
)

foreach ($r in $sheet3Rows) {
    $rowNum = $r[0]
    $ws3.Cells.Item($rowNum, 1).Style = "Normal"
    $ws3.Cells.Item($rowNum, 2).Value = $r[1]
    $ws3.Cells.Item($rowNum, 3).Value = $r[2]
    $ws3.Cells.Item($rowNum, 4).Value = $r[3]
    $ws3.Cells.Item($rowNum, 5).Value = $r[4]
    $ws3.Cells.Item($rowNum, 6).Value = $r[5]
    $ws3.Cells.Item($rowNum, 7).Value = $r[6]
}
